$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 7
$ws.Range("H1").Value = 8

$ws.Range("K24").Select()
